# Add columns I (I0) and J (IF) to Sheet1, matching header style of existing
# header row (B1:H1), and populate data rows 2-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy formatting from an existing header cell (H1) so the new
# header cells I1/J1 get the same bold/border/center style rather than a
# freshly synthesized one.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-12
$data = @{
    2  = @(9, 9)
    3  = @(11, 12)
    4  = @(7, 9)
    5  = @(8, 9)
    6  = @(9, 9)
    7  = @(7, 7)
    8  = @(4, 6)
    9  = @(5, 8)
    10 = @(6, 8)
    11 = @(4, 5)
    12 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
